$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Switched bumper pins from PB6/PB7 to PC6/PC7, with Left/Right pin numbers swapped.
$ws.Range("A22").Value = "C"
$ws.Range("B22").Value = 7
$ws.Range("B23").Value = 6

# Update the selection to match the post-edit active cell.
$ws.Range("F24").Select()
